$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 2 and row 3 for the columns that differ
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")
    $tmp = $cellRow2.Value2
    $cellRow2.Value2 = $cellRow3.Value2
    $cellRow3.Value2 = $tmp
}
